$d = $word.ActiveDocument

# First paragraph: "**ID__AFFARS_5343_topic_9__ID** "
$p1 = $d.Paragraphs(1)

# Add a paragraph border (top/left/bottom/right, space-only - no line drawn)
# and bump the left indent from 120 twips (6pt) to 225 twips (11.25pt),
# matching the pattern already used a couple of paragraphs below.
$pf = $p1.Range.ParagraphFormat
$pf.Borders.DistanceFromTop = 5
$pf.Borders.DistanceFromLeft = 5
$pf.Borders.DistanceFromBottom = 5
$pf.Borders.DistanceFromRight = 5
$pf.LeftIndent = 11.25

# Update the placeholder id text and drop the trailing stand-alone space run
# by including it in the search text but not in the replacement.
$d.Content.Find.Execute("**ID__AFFARS_5343_topic_9__ID** ", $true, $false, $false,
                         $false, $false, $false, 1, $false,
                         "**ID__AFFARS_5343_204_70_7__ID**", 2) | Out-Null
